# Add the new "Booking Management" team assignment (row 7, columns B and C)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B7").Value = "Booking Management"
$ws.Range("C7").Value = "Đặng Công Vinh + Võ Huy Tùng"

# Move the active selection to C7, matching the saved workbook state
$null = $ws.Range("C7").Select()
